# Added an excel process scope to make the program headless.
# Duplicate the "Total Sales" / "Profit Per Sale" columns (D:E) into new
# columns F:G on the Sales sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

$srcRange = $ws.Range("D1:E101")
$dstRange = $ws.Range("F1:G101")
$dstRange.Value = $srcRange.Value()
